# Sync attendance_reports: normalize "Recorded By" (column G) name/email
# ordering. The first two comma-separated entries of specific recorded-by
# strings are swapped (e.g. "System, dnasr281@gmail.com" becomes
# "dnasr281@gmail.com, System"), wherever that exact text appears in the
# "Session Analysis Results" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of exact old "Recorded By" text -> new text (first two entries swapped).
$map = @{
    "System, system, backup@backdoor.com" = "system, System, backup@backdoor.com"
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System"
    "backup@backdoor.com, System"         = "System, backup@backdoor.com"
    "admin@admin.com, System"             = "System, admin@admin.com"
    "admin@admin.com, dnasr281@gmail.com" = "dnasr281@gmail.com, admin@admin.com"
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count()

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value()
    if ($map.ContainsKey($current)) {
        $cell.Value = $map[$current]
    }
}
